# Apply the latest cryptocurrency market data scrape to Sheet1.
#
# Numeric-looking values (e.g. "1.009", "0.00001168", "1.200") are
# written with a leading apostrophe, which is exactly how Excel lets
# a user force text-entry for a value that would otherwise be
# auto-converted to a number. This preserves the exact textual
# formatting (leading/trailing zeros, "thousand.thousand.cents"
# style prices, etc.) that the source site renders, instead of
# losing it to Excel's numeric auto-detection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.530.83'

# Row 3
$ws.Range("D3").Value = '2.136.71'
$ws.Range("E3").Value = '  +1.66%  '

# Row 4
$ws.Range("D4").Value = '''1.009'
$ws.Range("E4").Value = '  +0.55%  '

# Row 5
$ws.Range("D5").Value = '''352.61'
$ws.Range("E5").Value = '  +5.03%  '

# Row 6
$ws.Range("E6").Value = '  +0.48%  '

# Row 7
$ws.Range("E7").Value = '  +0.58%  '

# Row 8
$ws.Range("D8").Value = '''0.4561'
$ws.Range("E8").Value = '  -0.02%  '

# Row 9
$ws.Range("D9").Value = '''53.53'
$ws.Range("E9").Value = '  -5.39%  '

# Row 10
$ws.Range("D10").Value = '''0.09147'
$ws.Range("E10").Value = '  +2.40%  '

# Row 11
$ws.Range("D11").Value = '''1.192'
$ws.Range("E11").Value = '  +1.13%  '

# Row 12
$ws.Range("D12").Value = '''25.51'
$ws.Range("E12").Value = '  +5.44%  '

# Row 13
$ws.Range("D13").Value = '2.135.50'
$ws.Range("E13").Value = '  +2.08%  '

# Row 14
$ws.Range("D14").Value = '''6.896'
$ws.Range("E14").Value = '  +1.04%  '

# Row 15
$ws.Range("D15").Value = '''8.210'
$ws.Range("E15").Value = '  +1.80%  '

# Row 16
$ws.Range("D16").Value = '''101.43'
$ws.Range("E16").Value = '  +4.34%  '

# Row 17
$ws.Range("D17").Value = '''0.00001168'
$ws.Range("E17").Value = '  +1.84%  '

# Row 18
$ws.Range("E18").Value = '  +0.50%  '

# Row 19
$ws.Range("D19").Value = '''0.06729'
$ws.Range("E19").Value = '  +1.43%  '

# Row 20
$ws.Range("D20").Value = '''20.48'
$ws.Range("E20").Value = '  +6.69%  '

# Row 21
$ws.Range("D21").Value = '''1.007'
$ws.Range("E21").Value = '  +0.48%  '

# Row 22
$ws.Range("D22").Value = '''6.376'
$ws.Range("E22").Value = '  +1.16%  '

# Row 23
$ws.Range("D23").Value = '30.634.79'
$ws.Range("E23").Value = '  +0.39%  '

# Row 24
$ws.Range("D24").Value = '''12.86'
$ws.Range("E24").Value = '  +4.08%  '

# Row 25
$ws.Range("D25").Value = '''2.378'
$ws.Range("E25").Value = '  +0.76%  '

# Row 26
$ws.Range("D26").Value = '2.391.44'
$ws.Range("E26").Value = '  +2.30%  '

# Row 27
$ws.Range("D27").Value = '''22.49'
$ws.Range("E27").Value = '  +1.25%  '

# Row 28
$ws.Range("D28").Value = '''2.601'
$ws.Range("E28").Value = '  +3.30%  '

# Row 29
$ws.Range("D29").Value = '''165.15'
$ws.Range("E29").Value = '  +1.59%  '

# Row 30
$ws.Range("D30").Value = '''135.78'
$ws.Range("E30").Value = '  +1.84%  '

# Row 31
$ws.Range("D31").Value = '''1.224'
$ws.Range("E31").Value = '  +1.19%  '

# Row 32
$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '''0.1082'
$ws.Range("E32").Value = '  +1.20%  '

# Row 33
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = '''1.709'
$ws.Range("E33").Value = '  +3.40%  '

# Row 34
$ws.Range("D34").Value = '''6.399'
$ws.Range("E34").Value = '  +0.27%  '

# Row 35
$ws.Range("D35").Value = '''4.032'
$ws.Range("E35").Value = '  +2.24%  '

# Row 36
$ws.Range("D36").Value = '''6.148'
$ws.Range("E36").Value = '  +4.41%  '

# Row 37
$ws.Range("D37").Value = '''10.43'
$ws.Range("E37").Value = '  +0.85%  '

# Row 38
$ws.Range("D38").Value = '''0.02643'
$ws.Range("E38").Value = '  +2.38%  '

# Row 39
$ws.Range("D39").Value = '''0.06984'
$ws.Range("E39").Value = '  +1.77%  '

# Row 40
$ws.Range("D40").Value = '''0.2359'
$ws.Range("E40").Value = '  +1.51%  '

# Row 41
$ws.Range("D41").Value = '''12.69'
$ws.Range("E41").Value = '  +0.13%  '

# Row 42
$ws.Range("D42").Value = '''0.6983'
$ws.Range("E42").Value = '  +1.57%  '

# Row 43
$ws.Range("D43").Value = '''1.277'
$ws.Range("E43").Value = '  +2.28%  '

# Row 44
$ws.Range("D44").Value = '''14.77'
$ws.Range("E44").Value = '  +5.06%  '

# Row 45
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '''0.6514'
$ws.Range("E45").Value = '  +1.93%  '

# Row 46
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = '''2.349'
$ws.Range("E46").Value = '  +0.78%  '

# Row 47
$ws.Range("D47").Value = '''0.00000000375'
$ws.Range("E47").Value = '  +10.46%  '

# Row 48
$ws.Range("D48").Value = '''3.717'
$ws.Range("E48").Value = '  +1.45%  '

# Row 49
$ws.Range("D49").Value = '''1.249'
$ws.Range("E49").Value = '  +0.16%  '

# Row 50
$ws.Range("D50").Value = '''83.81'
$ws.Range("E50").Value = '  +0.74%  '

# Row 51
$ws.Range("B51").Value = 'WEMIXTOKEN'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").Value = '''1.200'
$ws.Range("E51").Value = '  +0.23%  '
